$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-194 down to 71-195.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record (same shape as the
# surrounding rows - only D, M, N, O, P, R, S differ from the row that used
# to occupy position 70, the rest of the columns repeat the constant
# "Terminal Hortofrutícola Agro Chillán / Mango" context).
$ws.Range("A70").Value = 7
$ws.Range("B70").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C70").Value = "Ñuble"
$ws.Range("D70").Value = 45203
$ws.Range("E70").Value = 16
$ws.Range("F70").Value = "Fruta"
$ws.Range("G70").Value = 100108
$ws.Range("H70").Value = "Tropicales y subtropicales"
$ws.Range("I70").Value = 100108002
$ws.Range("J70").Value = "Mango"
$ws.Range("K70").Value = "Sin especificar"
$ws.Range("L70").Value = "Primera"
$ws.Range("M70").Value = 60
$ws.Range("N70").Value = 10000
$ws.Range("O70").Value = 10000
$ws.Range("P70").Value = 10000
$ws.Range("Q70").Value = "$/bandeja 4 kilos"
$ws.Range("R70").Value = "Brasil"
$ws.Range("S70").Value = 2500
$ws.Range("T70").Value = 4
